$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the row above (row 9) into the new row 10, then set the
# label value for column A to match the style used by the other year cells.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "2021年"

# Fill in the numeric data for the new 2021 row (columns B through R).
$ws.Range("B10").Value = 16674
$ws.Range("C10").Value = 597659
$ws.Range("D10").Value = 21389
$ws.Range("E10").Value = 27324
$ws.Range("F10").Value = 255406
$ws.Range("G10").Value = 5267
$ws.Range("H10").Value = 3371
$ws.Range("I10").Value = 48450
$ws.Range("J10").Value = 44524
$ws.Range("K10").Value = 106466
$ws.Range("L10").Value = 3794
$ws.Range("M10").Value = 7122
$ws.Range("N10").Value = 2883
$ws.Range("O10").Value = 7853
$ws.Range("P10").Value = 19152
$ws.Range("Q10").Value = 24727
$ws.Range("R10").Value = 3257
